# ADD results from server
# Update computed cost results for each year sheet (row 2 values)

$wb = $excel.ActiveWorkbook

# Index 0 corresponds to worksheet 1 ("2025"), index 1 to worksheet 2 ("2030"), etc.
$sheetData = @(
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 16189.88373611682; "G2" = 4231.516049511277; "I2" = 36871.39349763304; "M2" = 11723.69623729033; "N2" = 4440.772759108165; "O2" = 6938.835168893425 },
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 33088.99240739239; "G2" = 4231.516049511277; "I2" = 55489.39330640265; "M2" = 17557.90112633227; "N2" = 8420.406902678273; "O2" = 10344.8778666987 },
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 52205.12940273627; "G2" = 7864.0611328728;   "I2" = 71876.10593152183; "L2" = 9458.152972418126; "M2" = 24106.24159122616; "N2" = 11661.64784445888; "O2" = 13010.88149144449 },
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 52205.12940273627; "G2" = 7864.0611328728;   "I2" = 71876.10593152183; "L2" = 9458.152972418126; "M2" = 24106.24159122616; "N2" = 11661.64784445888; "O2" = 13010.88149144449 },
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 52205.12940273627; "G2" = 7864.0611328728;   "I2" = 71876.10593152183; "L2" = 9458.152972418126; "M2" = 24106.24159122616; "N2" = 11661.64784445888; "O2" = 13010.88149144449 },
    @{ "A2" = 4477.220219999998; "B2" = 6645.58835753044; "E2" = 52205.12940273627; "G2" = 7864.0611328728;   "I2" = 71876.10593152183; "L2" = 9458.152972418126; "M2" = 24106.24159122616; "N2" = 11661.64784445888; "O2" = 13010.88149144449 }
)

for ($i = 0; $i -lt $sheetData.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $cellValues = $sheetData[$i]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}
